$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-34: update Price (D) and/or Volume(1h) (E) values in place ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.812.11"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.808.35"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.31"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.35"
$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.806.48"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.16"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.449.32"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.748.98"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.59"
$ws.Range("E17").Value = "  +4.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.827.23"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.12"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.61"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("E22").Value = "  -5.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.702"
$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.59"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.11"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.01"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.957.61"
$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("E32").Value = "  +5.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.27"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.74"
$ws.Range("E34").Value = "  +0.09%  "

# --- Rows 35-51: coin list shifted up by one row, with a new coin (RenzoRestakedETH)
#     inserted and the last coin (EnergySwap) dropped ---
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.09"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.751.82"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.100"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +0.81%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.138"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.79"
$ws.Range("E42").Value = "  +1.21%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.14"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.87"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.299"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "150.47"
$ws.Range("E48").Value = "  +3.12%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.35"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "396.83"
$ws.Range("E50").Value = "  +2.79%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("E51").Value = "  -2.99%  "

